$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the post row for "「ママと一緒にお菓子作り」" (row 645).
# This shifts all subsequent rows up by one, so the sheet's used range
# shrinks from A1:C814 to A1:C813, matching the rest of the diff.
$ws.Rows.Item(645).EntireRow.Delete()
